$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9, pushing the existing rows 9-26 down to 10-27
# (this also carries the number formatting, e.g. the date style on column D,
# down from the row that used to be row 9)
$ws.Rows.Item(9).Insert()

# Populate the new row 9 with this week's data point
$ws.Range("A9").Value = 4
$ws.Range("B9").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C9").Value = "Los Lagos"
$ws.Range("D9").Value = 44897
$ws.Range("E9").Value = 10
$ws.Range("F9").Value = "Fruta"
$ws.Range("G9").Value = 100101
$ws.Range("H9").Value = "Berries"
$ws.Range("I9").Value = 100101001
$ws.Range("J9").Value = "Arándano (blue)"
$ws.Range("K9").Value = "Sin especificar"
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 400
$ws.Range("N9").Value = 6000
$ws.Range("O9").Value = 6500
$ws.Range("P9").Value = 6250
$ws.Range("Q9").Value = "$/bandeja 2 kilos"
$ws.Range("R9").Value = "Provincia de Curicó"
$ws.Range("S9").Value = 3125
$ws.Range("T9").Value = 2
